# Align WHO HPV workbook with current WHO guidance:
#  - extend single-dose series maxAge from 15y to 20y (WHO 2022 SAGE)
#  - change two-dose series threshold from >=15y to >=21y
#  - add a new three-dose "immunocompromised" risk series (obs 1022)

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell while forcing text storage even when
# the value looks like a bare integer (e.g. "1", "2", "3"), so it keeps
# the same shared-string / text semantics as the rest of the sheet.
function Set-TextValue {
    param($range, [string]$text)
    if ($text -match '^-?\d+$') {
        $range.NumberFormat = "@"
    }
    $range.Value = $text
}

# ---------------------------------------------------------------------
# Sheet "1-dose series": 15 years -> 20 years (maxAge for the single dose)
# ---------------------------------------------------------------------
$ws1d = $wb.Worksheets.Item("1-dose series")

Set-TextValue $ws1d.Cells.Item(7, 9) "20 years"   # I7 latestRecAge
Set-TextValue $ws1d.Cells.Item(9, 6) "20 years"   # F9 Dose 1 absMaxAge
$ws1d.Cells.Item(9, 7).Clear()                    # G9 (trailing n/a, dropped)
$ws1d.Cells.Item(9, 8).Clear()                    # H9 (trailing n/a, dropped)

# ---------------------------------------------------------------------
# Sheet "2-dose series": >=15y -> >=21y threshold
# ---------------------------------------------------------------------
$ws2d = $wb.Worksheets.Item("2-dose series")

Set-TextValue $ws2d.Cells.Item(1, 2) "WHO HPV 2-dose series (>=21y)"  # B1 Series Name
Set-TextValue $ws2d.Cells.Item(7, 8) "21 years"   # H7 absMinAge threshold
Set-TextValue $ws2d.Cells.Item(9, 3) "21 years"   # C9 Dose 1 absMinAge
Set-TextValue $ws2d.Cells.Item(9, 4) "21 years"   # D9 Dose 1 recMinAge
$ws2d.Cells.Item(9, 7).Clear()                    # G9 (trailing n/a, dropped)
$ws2d.Cells.Item(9, 8).Clear()                    # H9 (trailing n/a, dropped)

# Row 16 ("Age" row for Dose 2, "15 years + 5 months") is removed outright;
# the fixed age value is superseded by an interval-based rule, so everything
# below shifts up by one row.
$ws2d.Rows.Item(16).Delete()

# ---------------------------------------------------------------------
# New sheet: "3-dose series (immunocompromised)" (obs 1022)
# Excel worksheet names are capped at 31 characters, so the name is
# truncated to the maximum allowed length.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3dose = $wb.Worksheets.Add($null, $lastSheet)
$fullName = "3-dose series (immunocompromised)"
$ws3dose.Name = $fullName.Substring(0, [Math]::Min(31, $fullName.Length))

$dash = [char]0x2014

$rows = @(
    @{ R = 1;  Cells = @{ A = "Series Name"; B = "WHO HPV 3-dose series (immunocompromised)" } },
    @{ R = 2;  Cells = @{ A = "Target Disease"; B = "HPV" } },
    @{ R = 3;  Cells = @{ A = "Vaccine Group"; B = "HPV" } },
    @{ R = 4;  Cells = @{ A = "Series Type"; B = "Risk" } },
    @{ R = 5;  Cells = @{ A = "Equivalent Series Groups"; B = "n/a" } },
    @{ R = 6;  Cells = @{ A = "Gender"; B = "Female" } },
    @{ R = 7;  Cells = @{ A = "Select Patient Series"; B = "No"; C = "No"; D = "Immunocompromised"; E = "3"; F = "A"; G = "1"; H = "9 years"; I = "n/a" } },
    @{ R = 8;  Cells = @{ A = "Indication"; B = "Immunocompromised individual (1022)"; C = "Patient is immunocompromised $dash requires 3-dose HPV series"; D = "n/a"; E = "n/a" } },
    @{ R = 9;  Cells = @{ A = "Series Dose"; B = "Dose 1" } },
    @{ R = 10; Cells = @{ A = "Age"; B = "9 years"; C = "9 years"; D = "9 years"; E = "n/a"; F = "n/a" } },
    @{ R = 11; Cells = @{ A = "Preferable Vaccine"; B = "HPV, 9-valent (165)"; C = "9 years"; D = "n/a"; E = "n/a"; F = "n/a"; G = "Y" } },
    @{ R = 12; Cells = @{ A = "Allowable Vaccine"; B = "HPV, 9-valent (165)"; C = "9 years"; D = "n/a" } },
    @{ R = 13; Cells = @{ A = "Allowable Vaccine"; B = "HPV, quadrivalent (62)"; C = "9 years"; D = "n/a" } },
    @{ R = 14; Cells = @{ A = "Allowable Vaccine"; B = "HPV, bivalent (118)"; C = "9 years"; D = "n/a" } },
    @{ R = 15; Cells = @{ A = "Recurring Dose"; B = "No" } },
    @{ R = 16; Cells = @{ A = "Series Dose"; B = "Dose 2" } },
    @{ R = 17; Cells = @{ A = "Age"; B = "4 weeks"; C = "4 weeks"; D = "2 months"; E = "n/a" } },
    @{ R = 18; Cells = @{ A = "Preferable Vaccine"; B = "HPV, 9-valent (165)"; C = "9 years"; D = "n/a"; E = "n/a"; F = "n/a"; G = "Y" } },
    @{ R = 19; Cells = @{ A = "Allowable Vaccine"; B = "HPV, 9-valent (165)"; C = "9 years"; D = "n/a" } },
    @{ R = 20; Cells = @{ A = "Allowable Vaccine"; B = "HPV, quadrivalent (62)"; C = "9 years"; D = "n/a" } },
    @{ R = 21; Cells = @{ A = "Allowable Vaccine"; B = "HPV, bivalent (118)"; C = "9 years"; D = "n/a" } },
    @{ R = 22; Cells = @{ A = "Recurring Dose"; B = "No" } },
    @{ R = 23; Cells = @{ A = "Series Dose"; B = "Dose 3" } },
    @{ R = 24; Cells = @{ A = "Preferable Interval"; B = "Y"; C = "n/a"; D = "n/a"; E = "n/a"; F = "12 weeks"; G = "4 months"; H = "4 months"; I = "n/a" } },
    @{ R = 25; Cells = @{ A = "Preferable Vaccine"; B = "HPV, 9-valent (165)"; C = "9 years"; D = "n/a"; E = "n/a"; F = "n/a"; G = "Y" } },
    @{ R = 26; Cells = @{ A = "Allowable Vaccine"; B = "HPV, 9-valent (165)"; C = "9 years"; D = "n/a" } },
    @{ R = 27; Cells = @{ A = "Allowable Vaccine"; B = "HPV, quadrivalent (62)"; C = "9 years"; D = "n/a" } },
    @{ R = 28; Cells = @{ A = "Allowable Vaccine"; B = "HPV, bivalent (118)"; C = "9 years"; D = "n/a" } },
    @{ R = 29; Cells = @{ A = "Recurring Dose"; B = "No" } }
)

$colIndex = @{ A = 1; B = 2; C = 3; D = 4; E = 5; F = 6; G = 7; H = 8; I = 9 }

foreach ($rowDef in $rows) {
    $r = $rowDef.R
    foreach ($colLetter in @("A","B","C","D","E","F","G","H","I")) {
        if ($rowDef.Cells.ContainsKey($colLetter)) {
            $c = $colIndex[$colLetter]
            Set-TextValue $ws3dose.Cells.Item($r, $c) $rowDef.Cells[$colLetter]
        }
    }
}
